$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95, shifting existing rows 95:184 down to 96:185.
$ws.Rows("95:95").Insert()

# Populate the newly inserted row 95 with the new data record.
$ws.Range("A95").Value = 4
$ws.Range("B95").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C95").Value = "Los Lagos"
$ws.Range("D95").Value = 44566
$ws.Range("E95").Value = 10
$ws.Range("F95").Value = "Fruta"
$ws.Range("G95").Value = 100108
$ws.Range("H95").Value = "Tropicales y subtropicales"
$ws.Range("I95").Value = 100108005
$ws.Range("J95").Value = "Piña"
$ws.Range("K95").Value = "Caramelo"
$ws.Range("L95").Value = "Tercera"
$ws.Range("M95").Value = 20
$ws.Range("N95").Value = 19000
$ws.Range("O95").Value = 20000
$ws.Range("P95").Value = 19500
$ws.Range("Q95").Value = "$/caja 16 unidades"
$ws.Range("R95").Value = "Ecuador"
$ws.Range("S95").Value = 1219
$ws.Range("T95").Value = 16
